$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual data cell values that changed (star-count edits) ---
$ws.Range("I2").Value2 = "是"
$ws.Range("E3").Value2 = 5
$ws.Range("G4").Value2 = 5
$ws.Range("G5").Value2 = 4
$ws.Range("D6").Value2 = 5
$ws.Range("C7").Value2 = 1
$ws.Range("D9").Value2 = 5
$ws.Range("G9").Value2 = 3
$ws.Range("F10").Value2 = 3
$ws.Range("I12").Value2 = "是"
$ws.Range("D14").Value2 = 2
$ws.Range("G14").Value2 = 1
$ws.Range("D15").Value2 = 1
$ws.Range("G15").Value2 = 1

# --- Clear the 合计 (Total) column H contents, header included, keep formatting ---
$ws.Range("H1:H16").ClearContents()

# --- Hide column H since totals are no longer shown ---
$ws.Columns.Item(8).Hidden = $true

# --- Update the active selection to match the saved view state ---
$ws.Range("I3").Select()
